$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per diff
$ws.Range("J2").Value = 100
$ws.Range("C3").Value = 2
$ws.Range("J3").Value = 150

# Update selection to G12
$ws.Range("G12").Select()
